$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Ponderacion_nueva" (column C) values reflecting the adjusted
# rent (alquiler) weighting recalculation, plus a minor floating point
# cleanup for B27.

$ws.Range("C2").Value = 143.1428165540178
$ws.Range("C3").Value = 10.67712210810817
$ws.Range("C4").Value = 6.729731218203035
$ws.Range("C5").Value = 16.83502021585545
$ws.Range("C6").Value = 39.53410185804669
$ws.Range("C7").Value = 13.19493022187391
$ws.Range("C8").Value = 8.983007117436603
$ws.Range("C9").Value = 28.54096671521297
$ws.Range("C10").Value = 44.21014435684282
$ws.Range("C11").Value = 9.8573890481058
$ws.Range("C12").Value = 4.198458889925184
$ws.Range("C13").Value = 7.39422980500689
$ws.Range("C14").Value = 1.853309526961881
$ws.Range("C15").Value = 1.778860340836062
$ws.Range("C16").Value = 21.2172260332401
$ws.Range("C17").Value = 15.54324763531065
$ws.Range("C18").Value = 14.68945803271881
$ws.Range("C19").Value = 1.081889236679457
$ws.Range("C20").Value = 30.5035739835085
$ws.Range("C21").Value = 72.03592367767691
$ws.Range("C22").Value = 14.23642681501701
$ws.Range("C23").Value = 0.210675356483701
$ws.Range("C24").Value = 2.402174271485207
$ws.Range("C25").Value = 26.1950253396313
$ws.Range("C26").Value = 7.722123029007837
$ws.Range("B27").Value = 0.72
$ws.Range("C27").Value = 0.5702490852190403
$ws.Range("C28").Value = 10.88779746459187
$ws.Range("C29").Value = 24.97770194521232
$ws.Range("C30").Value = 10.77057959707462
$ws.Range("C31").Value = 13.19572223449227
$ws.Range("C32").Value = 3.690778801556566
$ws.Range("C33").Value = 1.816084933898972
$ws.Range("C34").Value = 5.258963785908927
$ws.Range("C35").Value = 2.42910270050944
$ws.Range("C36").Value = 89.70334915542848
$ws.Range("C37").Value = 8.016751723037675
$ws.Range("C38").Value = 25.97405381910892
$ws.Range("C39").Value = 5.280348126604641
$ws.Range("C40").Value = 3.508615899333817
$ws.Range("C41").Value = 12.20570646154254
$ws.Range("C42").Value = 0.8529975899734811
$ws.Range("C43").Value = 5.844261110876803
$ws.Range("C44").Value = 232.249104178437
